$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old row 4 (A4="admin", B4="password") content in place, leaving row 4 empty
$ws.Range("A4:B4").ClearContents()

# Add the two new rows at the bottom (row6 and row7 after the shift)
$ws.Range("A6").Value = "Admin"
$ws.Range("B6").Value = "Qedge123!@#"

$ws.Range("A7").Value = "admin"
$ws.Range("B7").Value = "password"

# Update the active selection to match the target view
$ws.Range("L11").Select()
